$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '27.906.00'
$ws.Range('E2').Value = '  +1.45%  '
$ws.Range('D3').Value = '1.638.60'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '213.59'
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('E6').Value = '  +0.91%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '23.58'
$ws.Range('E8').Value = '  +1.40%  '
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('E10').Value = '  +0.73%  '
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('D12').Value = '1.872.30'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('D13').Value = '1.638.95'
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.577'
$ws.Range('E14').Value = '  +4.76%  '
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '66.07'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '27.900.65'
$ws.Range('E17').Value = '  +1.48%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '232.05'
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.58'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.81'
$ws.Range('E22').Value = '  +3.42%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.36'
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('E24').Value = '  -3.67%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '151.75'
$ws.Range('E25').Value = '  +1.78%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.92'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '15.74'
$ws.Range('E27').Value = '  +1.56%  '
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0483'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.33'
$ws.Range('E32').Value = '  +2.01%  '
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.414.21'
$ws.Range('E33').Value = '  -3.72%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.10'
$ws.Range('E34').Value = '  +1.43%  '
$ws.Range('E35').Value = '  +1.61%  '
$ws.Range('E36').Value = '  +0.49%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.890'
$ws.Range('E37').Value = '  +1.92%  '
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.554'
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.912'
$ws.Range('E40').Value = '  -3.25%  '
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '66.35'
$ws.Range('E43').Value = '  -2.07%  '
$ws.Range('E44').Value = '  +3.77%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.44'
$ws.Range('E45').Value = '  +2.19%  '
$ws.Range('E46').Value = '  +0.32%  '
$ws.Range('D47').Value = '1.780.85'
$ws.Range('E47').Value = '  +1.16%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '88.15'
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('E49').Value = '  +1.02%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.61'
$ws.Range('E51').Value = '  -0.52%  '
